$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 ("grandes regiões e unidades da federação" header row).
# This shifts all rows below (7..37) up by one, so row 7 (norte, with its data)
# becomes the new row 6, etc. The orphaned shared string is cleaned up
# automatically since no cell references it any more.
$ws.Rows.Item(6).Delete()
